$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price column D, Volume(1h) column E)
# D-column numeric-looking values are written as text (matching the source
# inlineStr cells) by forcing a text NumberFormat before assignment so Excel
# does not auto-convert them into numbers (e.g. "1.00" -> 1).

$ws.Range("D2").Value = "70.447.63"
$ws.Range("E2").Value = "  -2.72%  "

$ws.Range("D3").Value = "2.520.45"
$ws.Range("E3").Value = "  -5.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.82"
$ws.Range("E5").Value = "  -3.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.85"
$ws.Range("E6").Value = "  -3.62%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  -2.52%  "

$ws.Range("D9").Value = "2.520.05"
$ws.Range("E9").Value = "  -5.10%  "

$ws.Range("E10").Value = "  -4.56%  "

$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("E12").Value = "  -3.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.83"
$ws.Range("E13").Value = "  -3.19%  "

$ws.Range("D14").Value = "2.978.48"
$ws.Range("E14").Value = "  -5.10%  "

$ws.Range("D15").Value = "70.309.38"
$ws.Range("E15").Value = "  -2.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("E16").Value = "  -3.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.98"
$ws.Range("E17").Value = "  -4.68%  "

$ws.Range("D18").Value = "2.511.24"
$ws.Range("E18").Value = "  -5.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.40"
$ws.Range("E19").Value = "  -7.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.67"
$ws.Range("E20").Value = "  +5.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.14"
$ws.Range("E21").Value = "  -4.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.95"
$ws.Range("E22").Value = "  -5.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.98"
$ws.Range("E23").Value = "  -4.25%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.10"
$ws.Range("E25").Value = "  -4.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.07"
$ws.Range("E26").Value = "  -5.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.13"
$ws.Range("E27").Value = "  -6.84%  "

$ws.Range("D28").Value = "2.649.25"
$ws.Range("E28").Value = "  -5.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").Value = "0.0₃0913"
$ws.Range("E30").Value = "  -5.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("E31").Value = "  -3.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "483.30"
$ws.Range("E32").Value = "  -2.39%  "

$ws.Range("E33").Value = "  -0.56%  "

$ws.Range("E34").Value = "  -3.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.39"
$ws.Range("E36").Value = "  -2.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("E37").Value = "  +3.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.91"
$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.59"
$ws.Range("E39").Value = "  -4.61%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  -3.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.74"
$ws.Range("E42").Value = "  -5.26%  "

$ws.Range("E43").Value = "  -6.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  -13.07%  "

$ws.Range("E45").Value = "  -7.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.29"
$ws.Range("E46").Value = "  -2.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.13"
$ws.Range("E47").Value = "  -8.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.53"
$ws.Range("E48").Value = "  -5.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.528"
$ws.Range("E49").Value = "  -5.09%  "

$ws.Range("E50").Value = "  -6.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.598"
$ws.Range("E51").Value = "  -1.14%  "
